$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Pekerjaan Ayah" / "Pekerjaan Ibu" / "Pekerjaan Wali" (father/mother/
# guardian occupation) columns are being dropped from the import template.
# They live at N (14), R (18) and V (22). Delete right-to-left so earlier
# column letters stay valid while we work.
$ws.Columns.Item(22).EntireColumn.Delete()
$ws.Columns.Item(18).EntireColumn.Delete()
$ws.Columns.Item(14).EntireColumn.Delete()

# After the deletions "Tanggal Lahir" (birth date) is column E. Rows 3 and 4
# should now hold their birth date as literal MM/DD/YYYY text instead of a
# real date serial, while keeping the existing date number format on the
# cell. Route the conversion through a scratch cell + Copy/PasteSpecial
# (values only) so the destination keeps its original style/number format
# but the pasted content is plain text, not a re-parsed date.
$ws.Range("Z1").Formula = "=TEXT(E3,""mm/dd/yyyy"")"
$ws.Range("Z1").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=TEXT(E4,""mm/dd/yyyy"")"
$ws.Range("Z1").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("Z1").Clear()

# Restore the active selection to where the author left off editing.
$ws.Range("S4").Select()
